# Power.xlsx / Sheet1 — fill in the "Down Power" (column D) formulas that
# propagate power draw up the supply tree, and point F7 at the newly
# added D7 instead of the shared SUM(G:S) formula ("getting on the DRC").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Formula = "=(B3-B4) * (D4  / B4)"
$ws.Range("D4").Formula = "=(B4 - B7) * F7 + (B4 - B10) * (D10 / B10) + (B4 - B13) * (D13 / B13) + (B4 - B14) * (D14 / B14)"
$ws.Range("D5").Formula = "=E5 + (B5 - B6) * F6"
$ws.Range("D6").Formula = "=E6"
$ws.Range("D7").Formula = "=(B7 - B8) * F8"
$ws.Range("F7").Formula = "=D7 / B7"
$ws.Range("D8").Formula = "=E8 + D9"
$ws.Range("D9").Formula = "=E9"
$ws.Range("D13").Formula = "=(B13 - B15)*F15 + (B13 - B17) * F17"
$ws.Range("D14").Formula = "=(B14 - B16) * F16"

$ws.Range("D3").Select() | Out-Null
